$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 1151
$ws.Range("J3").Value = 1235
$ws.Range("C4").Value = 1820
$ws.Range("I4").Value = 1754
$ws.Range("J4").Value = 266
$ws.Range("J6").Value = 1645
$ws.Range("C7").Value = 28363
$ws.Range("I7").Value = 26194
$ws.Range("J7").Value = 4390

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("J4").Value = 5
$ws.Range("J7").Value = 49

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("I4").Value = 30
$ws.Range("I7").Value = 449

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J3").Value = 67
$ws.Range("J6").Value = 54
$ws.Range("J7").Value = 162

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("J2").Value = 15
$ws.Range("J7").Value = 45

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J5").Value = 13
$ws.Range("J7").Value = 128
$ws.Range("J8").Value = 277
$ws.Range("J10").Value = 29
$ws.Range("J11").Value = 53
$ws.Range("J15").Value = 49
$ws.Range("J16").Value = 11
$ws.Range("J19").Value = 155
$ws.Range("J20").Value = 94
$ws.Range("J29").Value = 247
$ws.Range("J33").Value = 182
$ws.Range("J36").Value = 63
$ws.Range("J42").Value = 183
$ws.Range("J43").Value = 51
$ws.Range("J49").Value = 23
$ws.Range("J51").Value = 61
$ws.Range("J52").Value = 97
$ws.Range("J53").Value = 44
$ws.Range("J55").Value = 55
$ws.Range("C63").Value = 252
$ws.Range("I63").Value = 185
$ws.Range("J63").Value = 22
$ws.Range("J67").Value = 162
$ws.Range("J72").Value = 18
$ws.Range("J73").Value = 40
$ws.Range("J79").Value = 130
$ws.Range("J83").Value = 99
$ws.Range("J84").Value = 45
$ws.Range("J85").Value = 184
$ws.Range("J86").Value = 19
$ws.Range("J87").Value = 20
$ws.Range("J88").Value = 33
$ws.Range("J89").Value = 49
$ws.Range("J91").Value = 61
$ws.Range("J92").Value = 14
$ws.Range("J94").Value = 29
$ws.Range("J95").Value = 75
$ws.Range("J98").Value = 30
$ws.Range("I99").Value = 449
$ws.Range("J100").Value = 6
$ws.Range("C101").Value = 28363
$ws.Range("I101").Value = 26194
$ws.Range("J101").Value = 4390

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("J4").Value = 3
$ws.Range("J7").Value = 99

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("J4").Value = 3
$ws.Range("J7").Value = 75

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J2").Value = 43
$ws.Range("J3").Value = 49
$ws.Range("J7").Value = 182

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("J6").Value = 9
$ws.Range("J7").Value = 23

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J6").Value = 66
$ws.Range("J7").Value = 247

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J2").Value = 36
$ws.Range("J3").Value = 39
$ws.Range("J6").Value = 62
$ws.Range("J7").Value = 155

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J2").Value = 48
$ws.Range("J3").Value = 66
$ws.Range("J7").Value = 184

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J3").Value = 30
$ws.Range("J6").Value = 109
$ws.Range("J7").Value = 183

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("J6").Value = 11
$ws.Range("J7").Value = 29

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("J6").Value = 30
$ws.Range("J7").Value = 55

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("J3").Value = 26
$ws.Range("J6").Value = 12
$ws.Range("J7").Value = 61

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J2").Value = 36
$ws.Range("J3").Value = 44
$ws.Range("J7").Value = 130

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("J4").Value = 8
$ws.Range("J6").Value = 28
$ws.Range("J7").Value = 94

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("J6").Value = 29
$ws.Range("J7").Value = 63

$ws = $wb.Worksheets.Item("Wrigleyville")
$ws.Range("J3").Value = 2
$ws.Range("J6").Value = 6

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("J3").Value = 30
$ws.Range("J6").Value = 38
$ws.Range("J7").Value = 97

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("J2").Value = 6
$ws.Range("J6").Value = 16
$ws.Range("J7").Value = 29

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("J6").Value = 21
$ws.Range("J7").Value = 49

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("J6").Value = 15
$ws.Range("J7").Value = 30

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("J2").Value = 15
$ws.Range("J7").Value = 53

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("J2").Value = 16
$ws.Range("J7").Value = 40

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("J3").Value = 6
$ws.Range("J7").Value = 14

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("J3").Value = 11
$ws.Range("J7").Value = 33

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J2").Value = 93
$ws.Range("J3").Value = 90
$ws.Range("J7").Value = 277

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("J6").Value = 9
$ws.Range("J7").Value = 13

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("J2").Value = 4
$ws.Range("J3").Value = 5
$ws.Range("J7").Value = 19

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("J3").Value = 18
$ws.Range("J7").Value = 61

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("J3").Value = 9
$ws.Range("J7").Value = 51

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("J6").Value = 27
$ws.Range("J7").Value = 44

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("J6").Value = 4
$ws.Range("J7").Value = 18

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J2").Value = 43
$ws.Range("J3").Value = 45
$ws.Range("J7").Value = 128

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("J6").Value = 11
$ws.Range("J7").Value = 20

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("J6").Value = 7
$ws.Range("J7").Value = 11
